# The workbook's data table (barcode in column A, base value in column B)
# had its first 71 data rows (rows 2 through 72, i.e. the oldest 71
# records) removed. Every row below shifts up to fill the gap, the sheet's
# dimension/used-range shrinks from B472 to B401, and the now-unused
# shared-string barcode labels are dropped when Excel rewrites the
# sharedStrings table on save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2-72 (the top 71 data rows) - everything below shifts up.
$ws.Rows("2:72").Delete() | Out-Null

# Restore the active-cell selection to where it ended up after the
# deletion (matches the saved file's recorded selection).
$ws.Range("C11").Select() | Out-Null
